$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data appended to row 2: category, title, and study weight
$ws.Range("D2").Value = "student"
$ws.Range("E2").Value = "Reasons to Study in Poland"
$ws.Range("J2").Value = 20

# Row 2 height was recalculated slightly (auto height from wrapped text)
$ws.Rows.Item(2).RowHeight = 409.5

# User scrolled/selected cell L2 before saving
$ws.Range("L2").Select() | Out-Null
